# Auto-generated Excel COM-interop script applying the scheduled Sheets update.
# For each affected cell: numeric cells get a new .Value; cells removed entirely
# in the diff are cleared with .ClearContents(); cells newly introduced are set too.
$wb = $excel.ActiveWorkbook

# ---------- ALC ----------
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 19999
$ws.Range("I32").Value = 14997.8
$ws.Range("J32").Value = 23571.285
$ws.Range("K32").Value = 14997.8
$ws.Range("L32").Value = 23571.285
$ws.Range("M32").Value = -14671.8
$ws.Range("N32").Value = -24223.285
# Row 53
$ws.Range("H53").Value = 505.45456
$ws.Range("I53").Value = 231.75
$ws.Range("K53").Value = 231.75
$ws.Range("M53").Value = 405.25
# Row 64
$ws.Range("H64").Value = 4875
# Row 67
$ws.Range("H67").Value = 4875
# Row 74
$ws.Range("H74").Value = 3224.3
$ws.Range("I74").Value = 3382.5557
$ws.Range("K74").Value = 3382.5557
$ws.Range("M74").Value = -2446.5557
# Row 77
$ws.Range("H77").Value = 3224.3
$ws.Range("I77").Value = 3382.5557
$ws.Range("K77").Value = 16912.7785
$ws.Range("M77").Value = -12232.7785
# Row 86
$ws.Range("H86").Value = 7037.8335
$ws.Range("I86").Value = 6380.091
$ws.Range("K86").Value = 6380.091
$ws.Range("M86").Value = -5257.091
# Row 88
$ws.Range("H88").Value = 584027.9
$ws.Range("J88").Value = 907665.9
$ws.Range("L88").Value = 907665.9
$ws.Range("N88").Value = -908477.9
# Row 89
$ws.Range("H89").Value = 7037.8335
$ws.Range("I89").Value = 6380.091
$ws.Range("K89").Value = 31900.455
$ws.Range("M89").Value = -26284.455
# Row 91
$ws.Range("H91").Value = 584027.9
$ws.Range("J91").Value = 907665.9
$ws.Range("L91").Value = 907665.9
$ws.Range("N91").Value = -910473.9
# Row 100
$ws.Range("H100").Value = 3748.375
$ws.Range("I100").Value = 2972
$ws.Range("K100").Value = 2972
$ws.Range("M100").Value = -2431
# Row 106
$ws.Range("H106").Value = 2909.3333
$ws.Range("I106").Value = 3166.3333
$ws.Range("J106").Value = 2780.8333
$ws.Range("K106").Value = 3166.3333
$ws.Range("L106").Value = 2780.8333
$ws.Range("M106").Value = -2535.3333
$ws.Range("N106").Value = -4042.8333
# Row 113
$ws.Range("H113").Value = 100002216
$ws.Range("I113").Value = 33335032
$ws.Range("J113").Value = 200002990
$ws.Range("K113").Value = 33335032
$ws.Range("L113").Value = 200002990
$ws.Range("M113").Value = -33331778
$ws.Range("N113").Value = -200009498
# Row 120
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
# Row 137
$ws.Range("H137").Value = 5360.7896
$ws.Range("I137").Value = 3306.6667
$ws.Range("J137").Value = 8882.143
$ws.Range("K137").Value = 9920.000100000001
$ws.Range("L137").Value = 26646.429
$ws.Range("M137").Value = -7370.000100000001
$ws.Range("N137").Value = -31746.429

# ---------- ARM ----------
$ws = $wb.Worksheets.Item("ARM")
# Row 11
$ws.Range("H11").Value = 83337.336
$ws.Range("J11").Value = 83337.336
$ws.Range("L11").Value = 83337.336
$ws.Range("N11").Value = -83625.336
# Row 88
$ws.Range("H88").Value = 1652
$ws.Range("I88").Value = 1460.625
$ws.Range("K88").Value = 1460.625
$ws.Range("M88").Value = -1054.625
# Row 91
$ws.Range("H91").Value = 1652
$ws.Range("I91").Value = 1460.625
$ws.Range("K91").Value = 1460.625
$ws.Range("M91").Value = -56.625
# Row 102
$ws.Range("H102").Value = 12722.846
$ws.Range("I102").Value = 16599.889
$ws.Range("K102").Value = 16599.889
$ws.Range("M102").Value = -14977.889
# Row 125
$ws.Range("H125").Value = 46571.668
$ws.Range("J125").Value = 46571.668
$ws.Range("L125").Value = 46571.668
$ws.Range("N125").Value = -56411.668

# ---------- BSM ----------
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2200.25
$ws.Range("I94").Value = 1915.4286
$ws.Range("K94").Value = 1915.4286
$ws.Range("M94").Value = -1464.4286

# ---------- CRP ----------
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 866032.3
$ws.Range("I31").Value = 1479.2858
$ws.Range("K31").Value = 1479.2858
$ws.Range("M31").Value = -1184.2858
# Row 34
$ws.Range("H34").Value = 866032.3
$ws.Range("I34").Value = 1479.2858
$ws.Range("K34").Value = 1479.2858
$ws.Range("M34").Value = -1277.2858
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 99
$ws.Range("H99").Value = 3833
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
# Row 126
$ws.Range("H126").Value = 3833
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

# ---------- CUL ----------
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1577.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1577.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4732.5
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9072.5
# Row 114
$ws.Range("H114").Value = 1265.4166
$ws.Range("J114").Value = 928.5
$ws.Range("L114").Value = 2785.5
$ws.Range("N114").Value = -9293.5
# Row 131
$ws.Range("H131").Value = 3832.3845
$ws.Range("I131").Value = 1593
$ws.Range("J131").Value = 16149
$ws.Range("K131").Value = 4779
$ws.Range("L131").Value = 48447
$ws.Range("M131").Value = 261
$ws.Range("N131").Value = -58527

# ---------- GSM ----------
$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Range("H36").Value = 18201.3
$ws.Range("I36").Value = 13998.5
$ws.Range("K36").Value = 13998.5
$ws.Range("M36").Value = -13513.5
# Row 80
$ws.Range("H80").Value = 9847.272000000001
$ws.Range("I80").Value = 6092.5713
$ws.Range("J80").Value = 11599.467
$ws.Range("K80").Value = 6092.5713
$ws.Range("L80").Value = 11599.467
$ws.Range("M80").Value = -5094.5713
$ws.Range("N80").Value = -13595.467
# Row 83
$ws.Range("H83").Value = 9847.272000000001
$ws.Range("I83").Value = 6092.5713
$ws.Range("J83").Value = 11599.467
$ws.Range("K83").Value = 30462.8565
$ws.Range("L83").Value = 57997.33500000001
$ws.Range("M83").Value = -25470.8565
$ws.Range("N83").Value = -67981.33500000001
# Row 101
$ws.Range("H101").Value = 51832.43
$ws.Range("J101").Value = 51832.43
$ws.Range("L101").Value = 51832.43
$ws.Range("N101").Value = -58322.43
# Row 102
$ws.Range("H102").Value = 2931
$ws.Range("I102").Value = 1616.1
$ws.Range("J102").Value = 4809.4287
$ws.Range("K102").Value = 1616.1
$ws.Range("L102").Value = 4809.4287
$ws.Range("M102").Value = 5.900000000000091
$ws.Range("N102").Value = -8053.4287
# Row 105
$ws.Range("H105").Value = 115327.5
$ws.Range("J105").Value = 115327.5
$ws.Range("L105").Value = 115327.5
$ws.Range("N105").Value = -122315.5
# Row 106
$ws.Range("H106").Value = 113623
$ws.Range("J106").Value = 113623
$ws.Range("L106").Value = 113623
$ws.Range("N106").Value = -116147
# Row 126
$ws.Range("H126").Value = 9749.571
$ws.Range("I126").Value = 7356.2856
$ws.Range("J126").Value = 12142.857
$ws.Range("K126").Value = 22068.8568
$ws.Range("L126").Value = 36428.571
$ws.Range("M126").Value = -19598.8568
$ws.Range("N126").Value = -41368.571

# ---------- LTW ----------
$ws = $wb.Worksheets.Item("LTW")
# Row 19
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = 3
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 167
# Row 29
$ws.Range("H29").Value = 14999
$ws.Range("I29").Value = 14999
$ws.Range("K29").Value = 14999
$ws.Range("M29").Value = -14704
# Row 68
$ws.Range("H68").Value = 3982.5
$ws.Range("I68").Value = 1973.75
$ws.Range("K68").Value = 1973.75
$ws.Range("M68").Value = -1224.75
# Row 71
$ws.Range("H71").Value = 3982.5
$ws.Range("I71").Value = 1973.75
$ws.Range("K71").Value = 9868.75
$ws.Range("M71").Value = -6124.75
# Row 93
$ws.Range("H93").Value = 83335160
$ws.Range("I93").Value = 100001560
$ws.Range("K93").Value = 100001560
$ws.Range("M93").Value = -100000312
# Row 122
$ws.Range("H122").Value = 5288.2334
$ws.Range("I122").Value = 4802.227
$ws.Range("K122").Value = 14406.681
$ws.Range("M122").Value = -11956.681
# Row 132
$ws.Range("H132").Value = 369611.03
$ws.Range("I132").Value = 716963.9
$ws.Range("K132").Value = 2150891.7
$ws.Range("M132").Value = -2148361.7

# ---------- WVR ----------
$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 17000
$ws.Range("I32").Value = 17000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 17000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -16683
$ws.Range("N32").ClearContents()
# Row 64
$ws.Range("H64").Value = 64113.5
$ws.Range("J64").Value = 64113.5
$ws.Range("L64").Value = 64113.5
$ws.Range("N64").Value = -64609.5
# Row 67
$ws.Range("H67").Value = 64113.5
$ws.Range("J67").Value = 64113.5
$ws.Range("L67").Value = 64113.5
$ws.Range("N67").Value = -65829.5
